$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Cypher "StatQuery" text used in column C (rows 2-4 all show
# the same query), replacing the old breed-count query with the corrected
# one that also reports Program/Study/"Study Files" counts.
$newQuery = "MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)`nOPTIONAL MATCH (samp:sample)-->(c)`nOPTIONAL MATCH (diag:diagnosis)-->(c)`nOPTIONAL MATCH (f:file)-[*]->(c)`nOPTIONAL MATCH (sf:file)-->(s)`nWITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p`nWHERE demo.breed IN ['American Staffordshire Terrier']`nRETURN  `n    count(distinct p) AS Programs,`n    count(distinct s) AS Studies,`n    count(distinct c) AS Cases,`n    count(distinct samp) AS Samples,`n    count(distinct f) AS ``Case Files``,`n    count(distinct sf) AS ``Study Files``"
$ws.Range("C2").Value = $newQuery
$ws.Range("C3").Value = $newQuery
$ws.Range("C4").Value = $newQuery

# Update the active selection on the sheet (was B4, now B2) and drop the
# scrolled-away top-left cell so the view resets to show row 1.
$ws.Range("B2").Select()
